# Auto-generated edit script: updates cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep text-like numeric values (e.g. "583.82", "0.998") stored as text
# rather than being auto-coerced to numbers, matching the source data
# which stores prices/volumes as plain text.
$textRange1 = $ws.Range("D2:E51")
$textRange1.NumberFormat = "@"
$textRange2 = $ws.Range("B42:C45")
$textRange2.NumberFormat = "@"

$ws.Range("D2").Value = "60.529.39"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.622.19"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "583.82"
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").Value = "144.87"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "3.084.41"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "26.39"
$ws.Range("E14").Value = "  +14.00%  "
$ws.Range("D15").Value = "60.517.53"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "0.0000142"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "2.629.00"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "11.53"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "349.07"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "0.527"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "8.18"
$ws.Range("E27").Value = "  +7.84%  "
$ws.Range("E28").Value = "  +16.46%  "
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "6.58"
$ws.Range("E30").Value = "  +4.80%  "
$ws.Range("D31").Value = "169.89"
$ws.Range("E31").Value = "  +5.65%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "19.58"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "4.46"
$ws.Range("E34").Value = "  +6.41%  "
$ws.Range("D35").Value = "1.04"
$ws.Range("E35").Value = "  +9.36%  "
$ws.Range("D36").Value = "1.31"
$ws.Range("E36").Value = "  +8.99%  "
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").Value = "330.99"
$ws.Range("E38").Value = "  +13.26%  "
$ws.Range("D39").Value = "38.91"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("D41").Value = "0.873"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.18"
$ws.Range("E42").Value = "  +9.40%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "133.65"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.0999"
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "20.13"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "0.0557"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").Value = "0.610"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").Value = "20.41"
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "10.77"
$ws.Range("E51").Value = "  +0.79%  "

# Restore default (unstyled) cell style now that values are set as text
$textRange1.Style = "Normal"
$textRange2.Style = "Normal"

